$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "how are you"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = "Neutral"
